$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 3) {
        $shp = $candidate
        break
    }
}

$shp.Left = 169.7727
$shp.Top = 136.03976377952756
